# Edit: break out stock.yaml completed
# - backfill R1148/R1149 "backup" column from blank to 0
# - append 13 new weekly OHLCV rows (1150-1162)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R1148").Value = 0
$ws.Range("R1149").Value = 0

# Row 1150
$ws.Range("A1150").Value = 45474
$ws.Range("B1150").Value = 2796.10009765625
$ws.Range("C1150").Value = 2925
$ws.Range("D1150").Value = 2757.89990234375
$ws.Range("E1150").Value = 2891.10009765625
$ws.Range("F1150").Value = 2891.10009765625
$ws.Range("G1150").Value = 1634229
$ws.Range("H1150").Value = 2024
$ws.Range("I1150").Value = 7
$ws.Range("J1150").Value = 1
$ws.Range("K1150").Value = 0
$ws.Range("L1150").Value = 0
$ws.Range("M1150").Value = 0
$ws.Range("N1150").Value = 27
$ws.Range("O1150").Value = 0
$ws.Range("P1150").Value = 0
$ws.Range("Q1150").Value = 0
$ws.Range("A1150").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1151
$ws.Range("A1151").Value = 45481
$ws.Range("B1151").Value = 2891.85009765625
$ws.Range("C1151").Value = 2984
$ws.Range("D1151").Value = 2846.64990234375
$ws.Range("E1151").Value = 2948.449951171875
$ws.Range("F1151").Value = 2948.449951171875
$ws.Range("G1151").Value = 968158
$ws.Range("H1151").Value = 2024
$ws.Range("I1151").Value = 7
$ws.Range("J1151").Value = 8
$ws.Range("K1151").Value = 0
$ws.Range("L1151").Value = 0
$ws.Range("M1151").Value = 0
$ws.Range("N1151").Value = 28
$ws.Range("O1151").Value = 0
$ws.Range("P1151").Value = 0
$ws.Range("Q1151").Value = 0
$ws.Range("A1151").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1152
$ws.Range("A1152").Value = 45488
$ws.Range("B1152").Value = 2950
$ws.Range("C1152").Value = 3073.949951171875
$ws.Range("D1152").Value = 2935.10009765625
$ws.Range("E1152").Value = 2946.35009765625
$ws.Range("F1152").Value = 2946.35009765625
$ws.Range("G1152").Value = 1023274
$ws.Range("H1152").Value = 2024
$ws.Range("I1152").Value = 7
$ws.Range("J1152").Value = 15
$ws.Range("K1152").Value = 0
$ws.Range("L1152").Value = 0
$ws.Range("M1152").Value = 0
$ws.Range("N1152").Value = 29
$ws.Range("O1152").Value = 0
$ws.Range("P1152").Value = 0
$ws.Range("Q1152").Value = 0
$ws.Range("A1152").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1153
$ws.Range("A1153").Value = 45495
$ws.Range("B1153").Value = 2957.60009765625
$ws.Range("C1153").Value = 3257.14990234375
$ws.Range("D1153").Value = 2870
$ws.Range("E1153").Value = 3194.89990234375
$ws.Range("F1153").Value = 3194.89990234375
$ws.Range("G1153").Value = 2689112
$ws.Range("H1153").Value = 2024
$ws.Range("I1153").Value = 7
$ws.Range("J1153").Value = 22
$ws.Range("K1153").Value = 0
$ws.Range("L1153").Value = 0
$ws.Range("M1153").Value = 0
$ws.Range("N1153").Value = 30
$ws.Range("O1153").Value = 0
$ws.Range("P1153").Value = 0
$ws.Range("Q1153").Value = 0
$ws.Range("A1153").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1154
$ws.Range("A1154").Value = 45502
$ws.Range("B1154").Value = 3195
$ws.Range("C1154").Value = 3275.5
$ws.Range("D1154").Value = 3103
$ws.Range("E1154").Value = 3228.10009765625
$ws.Range("F1154").Value = 3228.10009765625
$ws.Range("G1154").Value = 2529052
$ws.Range("H1154").Value = 2024
$ws.Range("I1154").Value = 7
$ws.Range("J1154").Value = 29
$ws.Range("K1154").Value = 0
$ws.Range("L1154").Value = 0
$ws.Range("M1154").Value = 0
$ws.Range("N1154").Value = 31
$ws.Range("O1154").Value = 0
$ws.Range("P1154").Value = 0
$ws.Range("Q1154").Value = 0
$ws.Range("A1154").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1155
$ws.Range("A1155").Value = 45509
$ws.Range("B1155").Value = 3173.199951171875
$ws.Range("C1155").Value = 3372.85009765625
$ws.Range("D1155").Value = 3156.449951171875
$ws.Range("E1155").Value = 3338.39990234375
$ws.Range("F1155").Value = 3338.39990234375
$ws.Range("G1155").Value = 1588994
$ws.Range("H1155").Value = 2024
$ws.Range("I1155").Value = 8
$ws.Range("J1155").Value = 5
$ws.Range("K1155").Value = 0
$ws.Range("L1155").Value = 0
$ws.Range("M1155").Value = 0
$ws.Range("N1155").Value = 32
$ws.Range("O1155").Value = 0
$ws.Range("P1155").Value = 0
$ws.Range("Q1155").Value = 0
$ws.Range("A1155").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1156
$ws.Range("A1156").Value = 45516
$ws.Range("B1156").Value = 3338.39990234375
$ws.Range("C1156").Value = 3383.89990234375
$ws.Range("D1156").Value = 3260.89990234375
$ws.Range("E1156").Value = 3348.199951171875
$ws.Range("F1156").Value = 3348.199951171875
$ws.Range("G1156").Value = 1175746
$ws.Range("H1156").Value = 2024
$ws.Range("I1156").Value = 8
$ws.Range("J1156").Value = 12
$ws.Range("K1156").Value = 0
$ws.Range("L1156").Value = 0
$ws.Range("M1156").Value = 0
$ws.Range("N1156").Value = 33
$ws.Range("O1156").Value = 0
$ws.Range("P1156").Value = 0
$ws.Range("Q1156").Value = 0
$ws.Range("A1156").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1157
$ws.Range("A1157").Value = 45523
$ws.Range("B1157").Value = 3348.949951171875
$ws.Range("C1157").Value = 3391.949951171875
$ws.Range("D1157").Value = 3324.050048828125
$ws.Range("E1157").Value = 3349.35009765625
$ws.Range("F1157").Value = 3349.35009765625
$ws.Range("G1157").Value = 800957
$ws.Range("H1157").Value = 2024
$ws.Range("I1157").Value = 8
$ws.Range("J1157").Value = 19
$ws.Range("K1157").Value = 0
$ws.Range("L1157").Value = 0
$ws.Range("M1157").Value = 0
$ws.Range("N1157").Value = 34
$ws.Range("O1157").Value = 0
$ws.Range("P1157").Value = 0
$ws.Range("Q1157").Value = 0
$ws.Range("A1157").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1158
$ws.Range("A1158").Value = 45530
$ws.Range("B1158").Value = 3353.800048828125
$ws.Range("C1158").Value = 3574
$ws.Range("D1158").Value = 3311
$ws.Range("E1158").Value = 3485.14990234375
$ws.Range("F1158").Value = 3485.14990234375
$ws.Range("G1158").Value = 2116909
$ws.Range("H1158").Value = 2024
$ws.Range("I1158").Value = 8
$ws.Range("J1158").Value = 26
$ws.Range("K1158").Value = 0
$ws.Range("L1158").Value = 0
$ws.Range("M1158").Value = 0
$ws.Range("N1158").Value = 35
$ws.Range("O1158").Value = 0
$ws.Range("P1158").Value = 0
$ws.Range("Q1158").Value = 0
$ws.Range("A1158").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1159
$ws.Range("A1159").Value = 45537
$ws.Range("B1159").Value = 3503.449951171875
$ws.Range("C1159").Value = 3523
$ws.Range("D1159").Value = 3395
$ws.Range("E1159").Value = 3416.699951171875
$ws.Range("F1159").Value = 3416.699951171875
$ws.Range("G1159").Value = 1011821
$ws.Range("H1159").Value = 2024
$ws.Range("I1159").Value = 9
$ws.Range("J1159").Value = 2
$ws.Range("K1159").Value = 0
$ws.Range("L1159").Value = 0
$ws.Range("M1159").Value = 0
$ws.Range("N1159").Value = 36
$ws.Range("O1159").Value = 0
$ws.Range("P1159").Value = 0
$ws.Range("Q1159").Value = 0
$ws.Range("A1159").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1160
$ws.Range("A1160").Value = 45544
$ws.Range("B1160").Value = 3416
$ws.Range("C1160").Value = 3525
$ws.Range("D1160").Value = 3364.14990234375
$ws.Range("E1160").Value = 3451.699951171875
$ws.Range("F1160").Value = 3451.699951171875
$ws.Range("G1160").Value = 1328924
$ws.Range("H1160").Value = 2024
$ws.Range("I1160").Value = 9
$ws.Range("J1160").Value = 9
$ws.Range("K1160").Value = 0
$ws.Range("L1160").Value = 0
$ws.Range("M1160").Value = 0
$ws.Range("N1160").Value = 37
$ws.Range("O1160").Value = 0
$ws.Range("P1160").Value = 0
$ws.Range("Q1160").Value = 0
$ws.Range("A1160").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1161
$ws.Range("A1161").Value = 45551
$ws.Range("B1161").Value = 3460.39990234375
$ws.Range("C1161").Value = 3483.050048828125
$ws.Range("D1161").Value = 3320.949951171875
$ws.Range("E1161").Value = 3457.800048828125
$ws.Range("F1161").Value = 3457.800048828125
$ws.Range("G1161").Value = 1212520
$ws.Range("H1161").Value = 2024
$ws.Range("I1161").Value = 9
$ws.Range("J1161").Value = 16
$ws.Range("K1161").Value = 0
$ws.Range("L1161").Value = 0
$ws.Range("M1161").Value = 0
$ws.Range("N1161").Value = 38
$ws.Range("O1161").Value = 0
$ws.Range("P1161").Value = 0
$ws.Range("Q1161").Value = 0
$ws.Range("A1161").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1162
$ws.Range("A1162").Value = 45558
$ws.Range("B1162").Value = 3477.949951171875
$ws.Range("C1162").Value = 3520.35009765625
$ws.Range("D1162").Value = 3364.050048828125
$ws.Range("E1162").Value = 3482.60009765625
$ws.Range("F1162").Value = 3482.60009765625
$ws.Range("G1162").Value = 1408387
$ws.Range("H1162").Value = 2024
$ws.Range("I1162").Value = 9
$ws.Range("J1162").Value = 23
$ws.Range("K1162").Value = 0
$ws.Range("L1162").Value = 0
$ws.Range("M1162").Value = 0
$ws.Range("N1162").Value = 39
$ws.Range("O1162").Value = 0
$ws.Range("P1162").Value = 0
$ws.Range("Q1162").Value = 0
$ws.Range("A1162").NumberFormat = "YYYY-MM-DD HH:MM:SS"

